# Insert two new data rows at the top of the "Camote" price list block (row 372)
# for the week ending 44551 (1a nueva(o) / 2a nueva(o), origin Perú), pushing the
# existing rows 372-470 down to 374-472 and extending the sheet's used range to R472.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 372 (shifts 372:470 -> 374:472)
$ws.Rows.Item(372).Resize(2).Insert()

# New row 372: 1a nueva(o), fecha 44551, origen Perú
$ws.Cells.Item(372, 1).Value = 8
$ws.Cells.Item(372, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(372, 3).Value = "Coquimbo"
$ws.Cells.Item(372, 4).Value = 44551
$ws.Cells.Item(372, 5).Value = 4
$ws.Cells.Item(372, 6).Value = 100112045
$ws.Cells.Item(372, 7).Value = "Zapallo"
$ws.Cells.Item(372, 8).Value = "Camote"
$ws.Cells.Item(372, 9).Value = "1a nueva(o)"
$ws.Cells.Item(372, 10).Value = 1400
$ws.Cells.Item(372, 11).Value = 900
$ws.Cells.Item(372, 12).Value = 1000
$ws.Cells.Item(372, 13).Value = 950
$ws.Cells.Item(372, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(372, 15).Value = "Perú"
$ws.Cells.Item(372, 16).Value = 950
$ws.Cells.Item(372, 17).Value = 1
$ws.Cells.Item(372, 18).Value = "Hortaliza"

# New row 373: 2a nueva(o), fecha 44551, origen Perú
$ws.Cells.Item(373, 1).Value = 8
$ws.Cells.Item(373, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(373, 3).Value = "Coquimbo"
$ws.Cells.Item(373, 4).Value = 44551
$ws.Cells.Item(373, 5).Value = 4
$ws.Cells.Item(373, 6).Value = 100112045
$ws.Cells.Item(373, 7).Value = "Zapallo"
$ws.Cells.Item(373, 8).Value = "Camote"
$ws.Cells.Item(373, 9).Value = "2a nueva(o)"
$ws.Cells.Item(373, 10).Value = 800
$ws.Cells.Item(373, 11).Value = 800
$ws.Cells.Item(373, 12).Value = 850
$ws.Cells.Item(373, 13).Value = 825
$ws.Cells.Item(373, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(373, 15).Value = "Perú"
$ws.Cells.Item(373, 16).Value = 825
$ws.Cells.Item(373, 17).Value = 1
$ws.Cells.Item(373, 18).Value = "Hortaliza"
